# checklist_bank.xlsx - "updated 21 may 7pm"
#
# Marks additional checklist cells with "x" and re-centers a couple of
# previously vertical-only-aligned cells, then updates the saved
# selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (logger) and row 5 (exception): mark Customer columns (B:E) ---
$ws.Range("B4:E4").Value = "x"
$ws.Range("B5:E5").Value = "x"

# --- Transaction columns (M:N) for rows 6-11 need the same centered
#     alignment already used by the rest of the checklist grid ---
$ws.Range("M6:N11").HorizontalAlignment = -4108   # xlCenter

# Rows 6, 7, 10, 11 also get the "x" mark in the Transaction columns;
# rows 8 (code coverage) and 9 (jUnit) stay blank.
$ws.Range("M6:N7").Value = "x"
$ws.Range("M10:N11").Value = "x"

# --- Update the sheet's saved selection/scroll state ---
$ws.Range("N5").Select()

# --- Best-effort: resize the (host) window to match the author's saved
#     view. This mirrors the human action of resizing/maximizing the
#     Excel window; not all hosts persist this back into the workbook. ---
try {
    $excel.ActiveWindow.Width = 19380
    $excel.ActiveWindow.Height = 6075
} catch {
}
